$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quan Lộc có Thiên Lương thủ tọa tại Tý / Hợi
$ws.Range("A4").Value = "Quan Lộc có Thiên Lương thủ tọa tại Tý"
$ws.Range("A5").Value = "Quan Lộc có Thiên Lương thủ tọa tại Hợi"

$ws.Range("D4").Value = "Công việc nào cũng rất thích nhưng nhanh chán."
$ws.Range("C4").Value = "Hay thay đổi công việc."
$ws.Range("B4").Value = "Công việc có tính lưu động di chuyển."
$ws.Range("E4").Value = "Nên làm công việc có tính lưu động như du lịch, báo chí, hoặc lưu diễn..."

$ws.Range("D5").Value = "Công việc nào cũng rất thích nhưng nhanh chán."
$ws.Range("C5").Value = "Hay thay đổi công việc."
$ws.Range("B5").Value = "Công việc có tính lưu động di chuyển."
$ws.Range("E5").Value = "Nên làm công việc có tính lưu động như du lịch, báo chí, hoặc lưu diễn..."

# Quan Lộc có Thiên Đồng Thiên Lương đồng cung
$ws.Range("A6").Value = "Quan Lộc có Thiên Đồng Thiên Lương đồng cung"
$ws.Range("B6").Value = "Công việc liên quan đến y dược, chính trị, sư phạm."
$ws.Range("C6").Value = "Đều có danh tiếng trong ngành y dược, sư phạm."

# Quan Lộc có Tử Vi Thiên Tướng đồng cung
$ws.Range("A7").Value = "Quan Lộc có Tử Vi Thiên Tướng đồng cung`""
$ws.Range("B7").Value = "Công việc liên quan đến quân đội, cảnh sát. Có tiền đồ phát triển liên quan đến quân sự, quốc phòng."
$ws.Range("C7").Value = "Có tài lãnh binh, điều hành, điều khiển, lãnh đạo, chỉ đạo."

# Quan Lộc có Tham Lang Tử Vi đồng cung
$ws.Range("A8").Value = "Quan Lộc có Tham Lang Tử Vi đồng cung"
$ws.Range("B8").Value = "Công việc bình thường, nếu công việc có nhiều thành công rực rỡ dễ dính vào tai họa ẩn nấp"

$ws.Range("R16").Select()
